$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 86
$ws.Cells.Item($row, 1).Value = "2025-04-29 14:09:29"
$ws.Cells.Item($row, 2).Value = 253
